# Edited Chapter 3 up to 3.4.1
#
# - Renumber three hidden TOC bookmarks (Word doesn't let us assign
#   Bookmark.Name directly, so we re-add a bookmark with the new name
#   over the same range and drop the old one).
# - Refresh the cached PAGE-field text in the first section's header
#   and footer (18 -> 48, 17 -> 47).

$d = $word.ActiveDocument

function Rename-Bookmark($doc, $oldName, $newName) {
    $bm = $doc.Bookmarks($oldName)
    $r = $bm.Range
    $doc.Bookmarks.Add($newName, $r)
    $bm.Delete()
}

Rename-Bookmark $d "_Toc428458292" "_Toc430350701"
Rename-Bookmark $d "_Toc428458293" "_Toc430350702"
Rename-Bookmark $d "_Toc428458294" "_Toc430350703"

$sec = $d.Sections(1)

# header1.xml (default/primary header of section 1): page number 18 -> 48
$hdr = $sec.Headers(1)
$hdr.Range.Find.Execute("18", $false, $false, $false, $false, $false, $true, 1, $false, "48", 2) | Out-Null

# footer1.xml (first-page footer of section 1): page number 17 -> 47
$ftr = $sec.Footers(2)
$ftr.Range.Find.Execute("17", $false, $false, $false, $false, $false, $true, 1, $false, "47", 2) | Out-Null

Write-Host "Chapter 4 bookmarks renumbered and page numbers refreshed"
